# The "codeforiati:group-code" and "codeforiati:group-name" columns
# (columns C and D) were reordered in the source data, so every row's
# C and D values (including the header row) need to be swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal
    $ws.Cells.Item($r, 4).Value2 = $cVal
}
